$d = $word.ActiveDocument

# The title paragraph "List of questions and answers" is the first
# paragraph in the document. We insert a brand-new paragraph right
# after it (and before the "questions" bookmark / Heading2 section),
# carrying the "FirstParagraph" style and the new intro text.
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item(2)
$newPara.Range.Text = "If you would like practice questions and solutions to any of our topics, you can find them in our list of questions and answers."
$newPara.Style = "FirstParagraph"
